# "connected desktop app to backend"
# - Fix a typo in the rich-text description cell (F3): "restoran" -> "restaurant"
#   while preserving the existing run formatting (plain run, red "RestaurantModel"
#   run, and the trailing normal-colored run).
# - Update the desktop object-shape description (F4) to use proper
#   English/PascalCase field names now that the desktop app is wired up to
#   the backend.
# - Move the active selection to F5 (matches where the author ended up after
#   editing).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F3: fix " restoran" -> " restaurant", keeping rich-text runs intact ---
$cellF3 = $ws.Range("F3")
$before = $cellF3.Value2
$idx = $before.IndexOf(" restoran")
if ($idx -ge 0) {
    $oldRun = $cellF3.Characters($idx + 1, 9)   # " restoran" (space + 8 letters)
    $oldRun.Text = " restaurant"

    # Re-assert formatting on the (now longer) trailing run so the writer
    # keeps it as its own <r> instead of collapsing to plain text.
    $trailingRun = $cellF3.Characters($idx + 1, 11)  # " restaurant"
    $trailingRun.Font.Name = "Calibri"
    $trailingRun.Font.Size = 11
    $trailingRun.Font.Color = 0

    # Re-assert formatting on the red "RestaurantModel" run too, so it stays
    # a distinct run from the text around it.
    $fullAfter = $cellF3.Value2
    $modelIdx = $fullAfter.IndexOf("RestaurantModel")
    if ($modelIdx -ge 0) {
        $modelRun = $cellF3.Characters($modelIdx + 1, 15)  # "RestaurantModel"
        $modelRun.Font.Name = "Calibri"
        $modelRun.Font.Size = 11
        $modelRun.Font.Color = 255
    }
}

# --- F4: new object-shape description using proper names ---
$ws.Range("F4").Value = "string RestaurantName, string RestaurantDetails,string Food, string wine,doble Grade, string Image"

# --- move selection to F5 ---
$ws.Range("F5").Select()
